$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","G","H","I","J","L","N")
$data = @(
    @(2.66296437963365,0.2191847837540024,0.4359309698886165,0.1157484253495298,2.333094444815742,1.896143252811044,1.704586808734895,0.03259201963352965,0.6392722824210892,1.903177866949214),
    @(2.555238331714349,0.1959627064659628,0.4351887556438072,0.1162624875942235,2.322344464658812,1.898457203995122,1.70953589583798,0.03214406658179669,0.6309676664252208,1.926208370443058),
    @(2.490374876361443,0.1817686166681938,0.4349153733580664,0.1166100709673517,2.317169472202892,1.900815680177431,1.713471441166789,0.03186413073662209,0.6261763290399927,1.941036271554639),
    @(2.464264800847332,0.1760004084932234,0.4348498424295144,0.1167597601856318,2.315418145711604,1.902012305656655,1.715300336754034,0.03174881808051921,0.6243012043928218,1.947251390730739),
    @(2.459948717309828,0.1750435628868843,0.4348417325052196,0.116785102314533,2.315148905682022,1.902225222044791,1.715617609835633,0.03172959562437683,0.6239945153631368,1.948293828827854),
    @(2.490021441401211,0.1816907600549769,0.4349143038077585,0.1166120571326665,2.31714440684604,1.900830865010107,1.713495195260677,0.03186258060519975,0.6261507271204891,1.941119392056188),
    @(2.625554716407066,0.2111643341877993,0.4356372223036402,0.1159190510276282,2.329091403857774,1.896746318033649,1.706106895966691,0.03243857639249548,0.6363449780426578,1.910976097265702),
    @(2.901504321269215,0.2694840502115312,0.4385012184812638,0.1148130091162471,2.363877519090863,1.896190912105652,1.6987526035361,0.03352963724453062,0.6587794204554172,1.857322324484708),
    @(3.110479657482813,0.3126702124463066,0.4414874522650507,0.1141538948736383,2.396431129080014,1.900349989850895,1.697726602446352,0.03430832130875672,0.6767570713498259,1.821235928485342),
    @(3.206911534972846,0.3323949632543872,0.4430376509270815,0.113887232325057,2.412775933940651,1.903239005247258,1.698216622115623,0.03465770625957276,0.6852616121299064,1.805544423714574),
    @(3.243624772620194,0.3398758848741181,0.443652242917949,0.1137910126250503,2.41918737133787,1.904476761281558,1.698540250615835,0.03478932111286781,0.6885290632233989,1.799706801185794),
    @(3.235709178929142,0.3382642159252782,0.4435186536603339,0.1138115237382618,2.417796661128222,1.904203789023114,1.698464402920209,0.03476100606491883,0.6878232698592512,1.80095939001068),
    @(3.209928016061781,0.333010189775905,0.4430876613655954,0.1138792209417545,2.41329895067733,1.903337953058724,1.698240477299734,0.03466854807618347,0.6855294858555254,1.805062065039998),
    @(3.194161906270551,0.3297934621280376,0.4428272557743185,0.1139213069696403,2.410572920051749,1.902826335162985,1.698121311641145,0.0346118252750216,0.6841305941760254,1.807588676587038),
    @(3.104205122347594,0.3113827637002373,0.4413900008866705,0.1141719885032035,2.39539396300944,1.90018128072748,1.69771387243587,0.03428539131010311,0.6762078502392228,1.82227600234997),
    @(3.049369812274165,0.300108813004897,0.4405573965393188,0.1143342621149372,2.386476377887277,1.898814254414276,1.697709314461171,0.03408389914442544,0.6714311252453342,1.83147190365691),
    @(3.01795875631359,0.293631772096262,0.4400965529309104,0.11443072050829,2.381491725849344,1.89812179922211,1.697796717279971,0.03396754999474005,0.6687144014234008,1.836829309653082),
    @(3.00734563223557,0.2914400285344527,0.4399436192826016,0.1144639163456436,2.379828793060057,1.897903448476427,1.697841757370739,0.03392807761855465,0.6677998411754373,1.838654935616281),
    @(3.055193799108906,0.3013081730864258,0.4406441608131928,0.1143166646797606,2.387410705948128,1.898950063792341,1.697700478991806,0.03410539549414082,0.6719364359720856,1.830485927539685),
    @(3.217495234370119,0.3345531085038829,0.443213506185117,0.1138592075398641,2.414614002702479,1.90358836572878,1.698302498604917,0.03469572389201048,0.6862019509350574,1.803854174292576),
    @(3.324714372418043,0.3563482293512834,0.4450533917890596,0.1135879708362317,2.433687539267254,1.907457821761369,1.699500955321099,0.03507752097724648,0.6957990593580519,1.787057482901574),
    @(3.267384699185243,0.3447095057797469,0.4440567109749196,0.1137302003488543,2.423388796826202,1.905315811264359,1.698787495220856,0.03487411399329332,0.6906518419048666,1.795966414370994),
    @(3.052560418012604,0.3007659283850046,0.4406048791418868,0.1143246106206082,2.386987853420351,1.898888373254835,1.697704193126,0.03409567857944218,0.6717078931736324,1.830931467566717),
    @(2.825760354295085,0.2536490347576148,0.4375715617268838,0.1150852171858912,2.353243593978647,1.895540799088479,1.699975503799315,0.03323853958393386,0.6524481643832587,1.871252251818374)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $addr = "$($cols[$c])$rowNum"
        $ws.Range($addr).Value = $rowVals[$c]
    }
}
